$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.072082666666667
$ws.Range("H2").Value = 12.216248
$ws.Range("I2").Value = 0.3979924983064649
$ws.Range("J2").Value = 0.3979924983064649
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.719753333333333
$ws.Range("N2").Value = 8.15926
$ws.Range("O2").Value = 0.1068554218872926
$ws.Range("P2").Value = 0.1068554218872926
$ws.Range("Q2").Value = 11.07506040627556
$ws.Range("R2").Value = 99.67554365648
$ws.Range("S2").Value = 0.0425276563145149
$ws.Range("T2").Value = 0.0425276563145149
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.072082666666667
$ws.Range("H3").Value = 12.216248
$ws.Range("I3").Value = 0.3979924983064649
$ws.Range("J3").Value = 0.3979924983064649
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.856929333333333
$ws.Range("N3").Value = 17.570788
$ws.Range("O3").Value = 0.2301108145385953
$ws.Range("P3").Value = 0.2301108145385953
$ws.Range("Q3").Value = 23.84990041815822
$ws.Range("R3").Value = 214.649103763424
$ws.Range("S3").Value = 0.09158237796555113
$ws.Range("T3").Value = 0.09158237796555113
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.072082666666667
$ws.Range("H4").Value = 12.216248
$ws.Range("I4").Value = 0.3979924983064649
$ws.Range("J4").Value = 0.3979924983064649
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.639461666666667
$ws.Range("N4").Value = 19.918385
$ws.Range("O4").Value = 0.2608554492059968
$ws.Range("P4").Value = 0.2608554492059968
$ws.Range("Q4").Value = 27.03643676883111
$ws.Range("R4").Value = 243.32793091948
$ws.Range("S4").Value = 0.1038185119263498
$ws.Range("T4").Value = 0.1038185119263498
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.072082666666667
$ws.Range("H5").Value = 12.216248
$ws.Range("I5").Value = 0.3979924983064649
$ws.Range("J5").Value = 0.3979924983064649
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.531169
$ws.Range("N5").Value = 25.593507
$ws.Range("O5").Value = 0.3351780661555555
$ws.Range("P5").Value = 0.3351780661555555
$ws.Range("Q5").Value = 34.739625411304
$ws.Range("R5").Value = 312.656628701736
$ws.Range("S5").Value = 0.1333983559267791
$ws.Range("T5").Value = 0.1333983559267791
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.072082666666667
$ws.Range("H6").Value = 12.216248
$ws.Range("I6").Value = 0.3979924983064649
$ws.Range("J6").Value = 0.3979924983064649
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.705333666666667
$ws.Range("N6").Value = 5.116001
$ws.Range("O6").Value = 0.06700024821255984
$ws.Range("P6").Value = 0.06700024821255984
$ws.Range("Q6").Value = 6.944259664916444
$ws.Range("R6").Value = 62.498336984248
$ws.Range("S6").Value = 0.02666559617326995
$ws.Range("T6").Value = 0.02666559617326995
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.159473666666667
$ws.Range("H7").Value = 18.478421
$ws.Range("I7").Value = 0.6020075016935351
$ws.Range("J7").Value = 0.6020075016935351
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.719753333333333
$ws.Range("N7").Value = 8.15926
$ws.Range("O7").Value = 0.1068554218872926
$ws.Range("P7").Value = 0.1068554218872926
$ws.Range("Q7").Value = 16.75224903649556
$ws.Range("R7").Value = 150.77024132846
$ws.Range("S7").Value = 0.06432776557277774
$ws.Range("T7").Value = 0.06432776557277774
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.159473666666667
$ws.Range("H8").Value = 18.478421
$ws.Range("I8").Value = 0.6020075016935351
$ws.Range("J8").Value = 0.6020075016935351
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.856929333333333
$ws.Range("N8").Value = 17.570788
$ws.Range("O8").Value = 0.2301108145385953
$ws.Range("P8").Value = 0.2301108145385953
$ws.Range("Q8").Value = 36.07560199619422
$ws.Range("R8").Value = 324.680417965748
$ws.Range("S8").Value = 0.1385284365730441
$ws.Range("T8").Value = 0.1385284365730441
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.159473666666667
$ws.Range("H9").Value = 18.478421
$ws.Range("I9").Value = 0.6020075016935351
$ws.Range("J9").Value = 0.6020075016935351
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.639461666666667
$ws.Range("N9").Value = 19.918385
$ws.Range("O9").Value = 0.2608554492059968
$ws.Range("P9").Value = 0.2608554492059968
$ws.Range("Q9").Value = 40.89558929667611
$ws.Range("R9").Value = 368.060303670085
$ws.Range("S9").Value = 0.157036937279647
$ws.Range("T9").Value = 0.157036937279647
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.159473666666667
$ws.Range("H10").Value = 18.478421
$ws.Range("I10").Value = 0.6020075016935351
$ws.Range("J10").Value = 0.6020075016935351
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.531169
$ws.Range("N10").Value = 25.593507
$ws.Range("O10").Value = 0.3351780661555555
$ws.Range("P10").Value = 0.3351780661555555
$ws.Range("Q10").Value = 52.54751080138301
$ws.Range("R10").Value = 472.927597212447
$ws.Range("S10").Value = 0.2017797102287764
$ws.Range("T10").Value = 0.2017797102287764
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.159473666666667
$ws.Range("H11").Value = 18.478421
$ws.Range("I11").Value = 0.6020075016935351
$ws.Range("J11").Value = 0.6020075016935351
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.705333666666667
$ws.Range("N11").Value = 5.116001
$ws.Range("O11").Value = 0.06700024821255984
$ws.Range("P11").Value = 0.06700024821255984
$ws.Range("Q11").Value = 11.07506040627556
$ws.Range("R11").Value = 94.535620314421
$ws.Range("S11").Value = 0.0403346520392899
$ws.Range("T11").Value = 0.0403346520392899
